$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new task rows for the week (rows 26-27) ---
# Values are entered in the same order the original author typed them so the
# shared-string table grows in the same sequence as the target workbook.

# Row 26: "Base documents (agenda, notes and meeting)" worked on by Georgi,
# started 17.03. || 16:00, 40 minutes logged.
$ws.Range("A26").Value = "Base documents (agenda, notes and meeting)"
$ws.Range("B26").Value = "Georgi"
$ws.Range("C26").Value = "17.03. || 16:00"
$ws.Range("D26").Value = 40

# Row 27: "Setup document v1" worked on by Georgi and Ilia,
# started 21.03. || 17:00, 120 minutes logged.
$ws.Range("B27").Value = "Georgi and Ilia"
$ws.Range("C27").Value = "21.03. || 17:00"
$ws.Range("D27").Value = 120
$ws.Range("A27").Value = "Setup document v1"

# Give the two new rows the same look as the rest of the filled-in table
# (the blank template rows use a different style than populated ones).
$ws.Range("A4:D4").Copy()
$ws.Range("A26:D26").PasteSpecial(-4122)
$ws.Range("A23:D23").Copy()
$ws.Range("A27:D27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view so the newly entered rows are visible/selected ---
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("E27").Select()
